$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price / volume figures refreshed by the scheduled scraper run.
# Column D ("Price") cells that look like plain numbers must stay plain text
# (they use "." as a thousands separator and/or must preserve exact trailing
# zero digits), so we force the Text number format before writing them.

$ws.Range("D2").Value = "70.432.51"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Value = "3.621.00"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.04"
$ws.Range("E5").Value = "  -2.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.47"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("D7").Value = "3.614.74"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  -5.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.81"
$ws.Range("E11").Value = "  +16.74%  "

$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.33"
$ws.Range("E13").Value = "  -4.43%  "

$ws.Range("E14").Value = "  -2.29%  "

$ws.Range("D15").Value = "4.205.27"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "673.62"
$ws.Range("E16").Value = "  -4.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.91"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "3.625.85"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").Value = "70.456.97"
$ws.Range("E19").Value = "  -2.29%  "

$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.75"
$ws.Range("E21").Value = "  -4.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.44"
$ws.Range("E22").Value = "  -2.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.936"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.05"
$ws.Range("E24").Value = "  -4.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.58"
$ws.Range("E25").Value = "  -5.51%  "

$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("E27").Value = "  -2.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("E29").Value = "  -1.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.60"
$ws.Range("E30").Value = "  -3.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.04"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.54"
$ws.Range("E33").Value = "  +1.60%  "

$ws.Range("E34").Value = "  -6.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.98"
$ws.Range("E35").Value = "  -5.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "580.39"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("E37").Value = "  -3.16%  "

$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.18"
$ws.Range("E39").Value = "  -2.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "3.564.71"
$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0451"
$ws.Range("E42").Value = "  -0.43%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("E43").Value = "  -3.27%  "

$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.37"
$ws.Range("E45").Value = "  -4.68%  "

$ws.Range("E46").Value = "  -6.07%  "

$ws.Range("E47").Value = "  -4.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("E48").Value = "  +1.93%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.33"
$ws.Range("E50").Value = "  +2.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.91"
$ws.Range("E51").Value = "  -2.07%  "
